$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "New Method added at UI Drop Down": the SMOKE row's ExecutionFlag flips to
# Yes (it now runs) and the REGRESSION row's flips to No, i.e. A7/A8 swap.
$ws.Range("A7").Value = "Yes"
$ws.Range("A8").Value = "No"

# Match the widened columns (A/B) from the refreshed bestFit pass.
$ws.Columns.Item(1).ColumnWidth = 12.67
$ws.Columns.Item(2).ColumnWidth = 14.5

# The saved view now sits at 130% zoom with the cursor parked on A7.
$null = $ws.Range("A7").Select()
$excel.ActiveWindow.Zoom = 130
